$wb = $excel.ActiveWorkbook

# The "utilisateurs" sheet has a system account row (row 3) whose
# "ADRESSE MAIL" column (H) used to hold the placeholder "kcadmin".
# Update it to the new technical admin account name "ti_admin".
$wsUsers = $wb.Worksheets.Item("utilisateurs")
$wsUsers.Range("H3").Value = "ti_admin"

# Make "utilisateurs" the active sheet again (it was "zones"), with H3
# selected/active and scrolled so column H / row 1 is visible.
$wsUsers.Activate()
$excel.Goto($wsUsers.Range("H3"), $true)
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
